$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume symbol data pulled on Wed Jan  4 07:21:44 UTC 2023

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '254.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.38%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-4.59%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.320'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '3.25%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05850'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.80%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.712'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.80%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8663'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.62%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9113'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '5.48%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1429'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.83%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07169'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.16%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03181'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.80%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09221'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.71%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.54%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006075'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.96%'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.12%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.10%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.228'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.26%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3170'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.81%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03445'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.39%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1315'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.56%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.567'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.20%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04156'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.15%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.18%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.005039'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '21.66%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.001224'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.09%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001200'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '9.15%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '34.04%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03847'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.59%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1101'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.85%'
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003823'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-34.06%'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002380'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.01%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01097'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '26.88%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005237'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.01%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.04%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '54.86%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.14%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002098'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.04%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.04%'
